$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 74.609651
$ws.Range("H2").Value = 223.828953
$ws.Range("I2").Value = 0.1061386348809139
$ws.Range("J2").Value = 0.1061386348809139
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 658.1054191340806
$ws.Range("R2").Value = 5922.948772206726
$ws.Range("S2").Value = 0.006808949090969691
$ws.Range("T2").Value = 0.006808949090969693

# Row 3
$ws.Range("G3").Value = 74.609651
$ws.Range("H3").Value = 223.828953
$ws.Range("I3").Value = 0.1061386348809139
$ws.Range("J3").Value = 0.1061386348809139
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 4082.008532024285
$ws.Range("R3").Value = 36738.07678821857
$ws.Range("S3").Value = 0.04223364141269069
$ws.Range("T3").Value = 0.04223364141269071

# Row 4
$ws.Range("G4").Value = 74.609651
$ws.Range("H4").Value = 223.828953
$ws.Range("I4").Value = 0.1061386348809139
$ws.Range("J4").Value = 0.1061386348809139
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 1634.560420350997
$ws.Range("R4").Value = 14711.04378315897
$ws.Range("S4").Value = 0.01691163507349334
$ws.Range("T4").Value = 0.01691163507349335

# Row 5
$ws.Range("G5").Value = 74.609651
$ws.Range("H5").Value = 223.828953
$ws.Range("I5").Value = 0.1061386348809139
$ws.Range("J5").Value = 0.1061386348809139
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 3883.944082146209
$ws.Range("R5").Value = 34955.49673931588
$ws.Range("S5").Value = 0.04018440930376017
$ws.Range("T5").Value = 0.04018440930376018

# Row 6
$ws.Range("G6").Value = 597.374756
$ws.Range("I6").Value = 0.8498168837991085
$ws.Range("J6").Value = 0.8498168837991086
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 5269.232048512051
$ws.Range("R6").Value = 47423.08843660846
$ws.Range("S6").Value = 0.05451699944065468
$ws.Range("T6").Value = 0.05451699944065469

# Row 7
$ws.Range("G7").Value = 597.374756
$ws.Range("I7").Value = 0.8498168837991085
$ws.Range("J7").Value = 0.8498168837991086
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.3381507740050627
$ws.Range("T7").Value = 0.3381507740050628

# Row 8
$ws.Range("G8").Value = 597.374756
$ws.Range("I8").Value = 0.8498168837991085
$ws.Range("J8").Value = 0.8498168837991086
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 13087.38372565815
$ws.Range("R8").Value = 117786.4535309233
$ws.Range("S8").Value = 0.1354058588960445
$ws.Range("T8").Value = 0.1354058588960446

# Row 9
$ws.Range("G9").Value = 597.374756
$ws.Range("I9").Value = 0.8498168837991085
$ws.Range("J9").Value = 0.8498168837991086
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 31097.4534432514
$ws.Range("R9").Value = 279877.0809892626
$ws.Range("S9").Value = 0.3217432514573464
$ws.Range("T9").Value = 0.3217432514573466

# Row 10
$ws.Range("G10").Value = 30.48438
$ws.Range("H10").Value = 91.45313999999999
$ws.Range("I10").Value = 0.04336664808137267
$ws.Range("J10").Value = 0.04336664808137267
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 268.89196515532
$ws.Range("R10").Value = 2420.02768639788
$ws.Range("S10").Value = 0.002782034076124745
$ws.Range("T10").Value = 0.002782034076124745

# Row 11
$ws.Range("G11").Value = 30.48438
$ws.Range("H11").Value = 91.45313999999999
$ws.Range("I11").Value = 0.04336664808137267
$ws.Range("J11").Value = 0.04336664808137267
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 1667.8472233233
$ws.Range("R11").Value = 15010.6250099097
$ws.Range("S11").Value = 0.01725602996867255
$ws.Range("T11").Value = 0.01725602996867255

# Row 12
$ws.Range("G12").Value = 30.48438
$ws.Range("H12").Value = 91.45313999999999
$ws.Range("I12").Value = 0.04336664808137267
$ws.Range("J12").Value = 0.04336664808137267
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 667.8567761554
$ws.Range("R12").Value = 6010.710985398599
$ws.Range("S12").Value = 0.006909839452294167
$ws.Range("T12").Value = 0.006909839452294167

# Row 13
$ws.Range("G13").Value = 30.48438
$ws.Range("H13").Value = 91.45313999999999
$ws.Range("I13").Value = 0.04336664808137267
$ws.Range("J13").Value = 0.04336664808137267
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 1586.92107136242
$ws.Range("R13").Value = 14282.28964226178
$ws.Range("S13").Value = 0.0164187445842812
$ws.Range("T13").Value = 0.0164187445842812

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4764796666666666
$ws.Range("H14").Value = 1.429439
$ws.Range("I14").Value = 0.0006778332386049212
$ws.Range("J14").Value = 0.0006778332386049213
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 4.202859101170889
$ws.Range("R14").Value = 37.825731910538
$ws.Range("S14").Value = 0.00004348399636952519
$ws.Range("T14").Value = 0.00004348399636952519

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4764796666666666
$ws.Range("H15").Value = 1.429439
$ws.Range("I15").Value = 0.0006778332386049212
$ws.Range("J15").Value = 0.0006778332386049213
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 26.06893395962166
$ws.Range("R15").Value = 234.620405636595
$ws.Range("S15").Value = 0.0002697167338638052
$ws.Range("T15").Value = 0.0002697167338638053

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4764796666666666
$ws.Range("H16").Value = 1.429439
$ws.Range("I16").Value = 0.0006778332386049212
$ws.Range("J16").Value = 0.0006778332386049213
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 10.43879436234556
$ws.Range("R16").Value = 93.94914926110999
$ws.Range("S16").Value = 0.000108002786966614
$ws.Range("T16").Value = 0.000108002786966614

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4764796666666666
$ws.Range("H17").Value = 1.429439
$ws.Range("I17").Value = 0.0006778332386049212
$ws.Range("J17").Value = 0.0006778332386049213
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 24.80403482403366
$ws.Range("R17").Value = 223.2363134163029
$ws.Range("S17").Value = 0.0002566297214049767
$ws.Range("T17").Value = 0.0002566297214049768
